$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update summary figures -------------------------------------------------
# Valor Mora total
$ws.Range("E11").Value = 239148
# Cant. Periodos count (one more period added)
$ws.Range("F13").Value = 5

# --- Insert a new data row for the new period (2508) ------------------------
# The table currently ends at row 19 (the row with the closing/bottom border).
# Insert a fresh row at 20 so the old row 19 (bottom-border styled) shifts down
# and a blank row 20 appears in its place; the footer rows below move down too.
$ws.Rows(20).Insert()

# Give the new row 20 the same "closing" (bottom border) formatting that row 19
# used to have, by copying formats only from the row that is still row 19.
$ws.Range("B19:J19").Copy()
$ws.Range("B20:J20").PasteSpecial(-4122)

# Re-style row 19 to look like a normal "middle" data row (same as rows 16-18)
# instead of the closing row, since it's no longer the last period.
$ws.Range("B16:J16").Copy()
$ws.Range("B19:J19").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Reorder the period rows into ascending order and add the new period ----
# Old order was 2507, 2506, 2505, 2504 (rows 16-19); new order is
# 2504, 2505, 2506, 2507, 2508 (rows 16-20).
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1041973496"
$ws.Range("D16").Value = "JUAN DANIEL CEBALLO SIOLO"
$ws.Range("E16").Value = "2504"
$ws.Range("F16").Value = 11388
$ws.Range("G16").Value = 1423500

$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1041973496"
$ws.Range("D17").Value = "JUAN DANIEL CEBALLO SIOLO"
$ws.Range("E17").Value = "2505"
$ws.Range("F17").Value = 56940
$ws.Range("G17").Value = 1423500

$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1041973496"
$ws.Range("D18").Value = "JUAN DANIEL CEBALLO SIOLO"
$ws.Range("E18").Value = "2506"
$ws.Range("F18").Value = 56940
$ws.Range("G18").Value = 1423500

$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1041973496"
$ws.Range("D19").Value = "JUAN DANIEL CEBALLO SIOLO"
$ws.Range("E19").Value = "2507"
$ws.Range("F19").Value = 56940
$ws.Range("G19").Value = 1423500

$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "1041973496"
$ws.Range("D20").Value = "JUAN DANIEL CEBALLO SIOLO"
$ws.Range("E20").Value = "2508"
$ws.Range("F20").Value = 56940
$ws.Range("G20").Value = 1423500
